$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: "Product"
#   - keep header (ProductName) + the "RAW_" row (was row 4),
#     drop "QSP_", "ISI_", "CIA_"
#   - add a second column "ExpectedError" to the table
#   - move selection to A3
# ============================================================
$ws1 = $wb.Worksheets.Item(1)

# Row2 (QSP_) and Row3 (ISI_) go away; "RAW_" (row4) shifts up to row2.
$ws1.Range("A2:A3").EntireRow.Delete() | Out-Null
# What is now row3 ("CIA_") also goes away.
$ws1.Range("A3").EntireRow.Delete() | Out-Null

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:B2"))
$ws1.Range("B1").Value = "ExpectedError"

$ws1.Range("A3").Select() | Out-Null

# ============================================================
# Sheet 2: "Organization"
#   - keep header (OrgName) + the "hamas_" row (was row 4),
#     drop "ISIS_", "Mujhahideen_"
#   - shrink the conditional formatting to the single remaining row
#   - becomes the active tab, selection A2:A3 (active cell A2)
# ============================================================
$ws2 = $wb.Worksheets.Item(2)

# Row2 (ISIS_) and Row3 (Mujhahideen_) go away; "hamas_" (row4) shifts to row2.
$ws2.Range("A2:A3").EntireRow.Delete() | Out-Null

$fc2 = $ws2.Range("A2:A4").FormatConditions
if ($fc2.Count -gt 0) {
    $fc2.Item(1).ModifyAppliesToRange($ws2.Range("A2"))
}

$ws2.Range("A2:A3").Select() | Out-Null

# ============================================================
# Sheet 3: "Contact" -> "Opportunities"
#   - brand new data: Opportunity Name / AsiaCup
#   - drop the second column entirely
#   - table shrinks to a single column but keeps ref ending at row 3
#   - no longer the active tab, selection B4
# ============================================================
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Opportunities"

$ws3.Range("B1:B3").ClearContents() | Out-Null
$ws3.Range("A3").ClearContents() | Out-Null

$ws3.Range("A1").Value = "Opportunity Name"
$ws3.Range("A2").Value = "AsiaCup"

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:A3"))

$ws3.Range("B4").Select() | Out-Null

# ============================================================
# Workbook level: Organization becomes the active sheet/tab
# ============================================================
$ws2.Activate() | Out-Null
